$d = $word.ActiveDocument

$replacements = @(
    @("15×30=", "19×65="),
    @("42×35=", "12×15="),
    @("39×23=", "96×99="),
    @("48×20=", "85×97="),
    @("41×81=", "49×94="),
    @("52×37=", "84×86="),
    @("39×21=", "88×91="),
    @("92×16=", "58×71="),
    @("57×85=", "19×87="),
    @("91×69=", "83×88="),
    @("38×24=", "93×49="),
    @("21×84=", "77×36="),
    @("20×88=", "34×40="),
    @("28×63=", "93×15="),
    @("43×98=", "13×71="),
    @("79×25=", "76×53="),
    @("45×60=", "54×87="),
    @("42×27=", "16×46="),
    @("63×92=", "91×50="),
    @("91×63=", "90×56="),
    @("33×23=", "76×63="),
    @("15×24=", "53×85="),
    @("76×48=", "30×98="),
    @("18×37=", "11×56="),
    @("52×79=", "95×66=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
